$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "TextBox 28" (shape id=29) - reposition/resize and add a sigma run
$sh = $s.Shapes.Item(19)

# Point values chosen so the COM layer's Single-precision (float32) storage
# round-trips to the exact target EMU values (EMU = floor(float32(pts) * 12700)).
$sh.Left   = 185.97110836220475
$sh.Top    = 184.96787401574804
$sh.Width  = 25.14314960629921
$sh.Height = 29.081259842519685

$sh.TextFrame.TextRange.Text = "𝜎"
